$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the "Price" (D) and "Volume(1h)" (E) columns for each crypto row
# with the latest scraped figures. Price values are entered with a leading
# apostrophe so numeric-looking strings (e.g. "0.9991") are stored as text,
# matching the sheet's existing inline-string cell format instead of being
# auto-coerced to a Number type.
$ws.Range("D2").Value = "'28.653.16"
$ws.Range("E2").Value = "  +2.65%  "
$ws.Range("D3").Value = "'1.915.77"
$ws.Range("E3").Value = "  +5.77%  "
$ws.Range("D4").Value = "'0.9991"
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").Value = "'313.82"
$ws.Range("E5").Value = "  +1.36%  "
$ws.Range("D6").Value = "'0.9991"
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").Value = "'0.5051"
$ws.Range("E7").Value = "  +2.26%  "
$ws.Range("D8").Value = "'0.3972"
$ws.Range("E8").Value = "  +2.68%  "
$ws.Range("D9").Value = "'0.09726"
$ws.Range("E9").Value = "  -0.81%  "
$ws.Range("D10").Value = "'1.165"
$ws.Range("E10").Value = "  +5.82%  "
$ws.Range("D11").Value = "'41.68"
$ws.Range("E11").Value = "  +2.01%  "
$ws.Range("D12").Value = "'6.579"
$ws.Range("E12").Value = "  +2.25%  "
$ws.Range("D13").Value = "'21.25"
$ws.Range("E13").Value = "  +3.62%  "
$ws.Range("D14").Value = "'1.923.94"
$ws.Range("E14").Value = "  +6.39%  "
$ws.Range("D15").Value = "'7.581"
$ws.Range("E15").Value = "  +4.03%  "
$ws.Range("D16").Value = "'0.9999"
$ws.Range("E16").Value = "  -0.08%  "
$ws.Range("E17").Value = "  +0.28%  "
$ws.Range("D18").Value = "'94.07"
$ws.Range("E18").Value = "  +1.60%  "
$ws.Range("D19").Value = "'0.06627"
$ws.Range("E19").Value = "  +0.39%  "
$ws.Range("D20").Value = "'18.08"
$ws.Range("E20").Value = "  +5.91%  "
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("D22").Value = "'6.285"
$ws.Range("E22").Value = "  +5.96%  "
$ws.Range("D23").Value = "'28.706.71"
$ws.Range("E23").Value = "  +2.63%  "
$ws.Range("D24").Value = "'11.47"
$ws.Range("E24").Value = "  +2.97%  "
$ws.Range("D25").Value = "'2.280"
$ws.Range("D26").Value = "'2.777"
$ws.Range("E26").Value = "  +16.18%  "
$ws.Range("D27").Value = "'2.134.61"
$ws.Range("E27").Value = "  +5.67%  "
$ws.Range("D28").Value = "'21.46"
$ws.Range("E28").Value = "  +4.36%  "
$ws.Range("D29").Value = "'159.29"
$ws.Range("E29").Value = "  +0.31%  "
$ws.Range("D30").Value = "'128.73"
$ws.Range("E30").Value = "  +1.14%  "
$ws.Range("D31").Value = "'1.120"
$ws.Range("E31").Value = "  +7.60%  "
$ws.Range("D32").Value = "'0.1078"
$ws.Range("E32").Value = "  +1.75%  "
$ws.Range("D33").Value = "'5.737"
$ws.Range("E33").Value = "  +2.87%  "
$ws.Range("D34").Value = "'3.635"
$ws.Range("E34").Value = "  +0.13%  "
$ws.Range("D35").Value = "'9.835"
$ws.Range("E35").Value = "  +9.01%  "
$ws.Range("D36").Value = "'0.06818"
$ws.Range("E36").Value = "  +0.91%  "
$ws.Range("E37").Value = "  +5.13%  "
$ws.Range("D38").Value = "'0.2224"
$ws.Range("E38").Value = "  +4.32%  "
$ws.Range("D39").Value = "'5.126"
$ws.Range("E39").Value = "  +3.82%  "
$ws.Range("E40").Value = "  +3.57%  "
$ws.Range("D41").Value = "'0.6441"
$ws.Range("E41").Value = "  +3.99%  "
$ws.Range("D42").Value = "'1.199"
$ws.Range("E42").Value = "  +5.00%  "
$ws.Range("D43").Value = "'0.9988"
$ws.Range("D44").Value = "'13.85"
$ws.Range("E44").Value = "  +6.24%  "
$ws.Range("D45").Value = "'0.6112"
$ws.Range("E45").Value = "  +4.21%  "
$ws.Range("D46").Value = "'1.286"
$ws.Range("E46").Value = "  +0.62%  "
$ws.Range("D47").Value = "'3.656"
$ws.Range("E47").Value = "  -0.86%  "
$ws.Range("D48").Value = "'2.050"
$ws.Range("E48").Value = "  +6.24%  "
$ws.Range("D49").Value = "'125.01"
$ws.Range("E49").Value = "  +2.16%  "
$ws.Range("E50").Value = "  +3.23%  "
$ws.Range("D51").Value = "'78.51"
$ws.Range("E51").Value = "  +6.70%  "
